# Replace curly quotes with straight single quotes in English (en_US) lines
# for act13d5_05_beg story sheet (column C), per commit 20210731.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = "[name=""Spokesman Czarny""]  You are a good person, and they would deprive you of your own living to crown their ancient sense of 'honor.'`n"
$ws.Range("C26").Value = "[name=""Spokesman Czarny""]  And so those 'bad people' would ignore that, their country and all its citizens, for their own inconsequential, self-centered consolation.`n"
$ws.Range("C37").Value = "[name=""Old Craftsman""]  Hey, Marcin. Think about carving it into this seat here, 'Kowal ONLY.' Don’t wanna wait for a spot next time I come by.`n"
$ws.Range("C44").Value = "[name=""Tourist""]  Hold on... you’re Zofia! The 'Whislash' Zofia!  `n"
$ws.Range("C45").Value = "[name=""Tourist""]  Hey, the rumors were true! 'Whislash' Zofia is the actual coach of the new Nearl! `n"
$ws.Range("C59").Value = "[name=""Tourist""]  Marcin...? You’re the real 'Tremoriron' Marcin... *sobbing*... Today’s just too blessed... Give me a 'Marcin Special'... `n"
$ws.Range("C103").Value = "[name=""Maria""]  ...'Flametail' said so too. I haven’t seen what the significance of the Major really is...`n"
$ws.Range("C108").Value = "[name=""Zofia""]  That’s not a 'sacrifice' at all. And I don’t want to watch you suffer for any of this, either...`n"
$ws.Range("C123").Value = "[name=""Bald Marcin""]  The corps were only interested in giving 'the new Nearl knight' the chance to take root and start growing. When they decide to reap you, they’ll come out with their sickles, no second thoughts.`n"

